# "Add files via upload" -- populate the previously-empty Sheet1 with the
# "Data Fields Details" planning table (Table/Columns layout for the
# User, Group and Chat entities).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
# Written in the exact order their text first appears so the generated
# shared-strings table comes out in the same order as the target file.
$ws.Range("B2").Value  = '"User"'
$ws.Range("B4").Value  = "Username"
$ws.Range("B5").Value  = "Password"
$ws.Range("B6").Value  = "Email"
$ws.Range("B8").Value  = "Group"
$ws.Range("B3").Value  = "Unique_UID"
$ws.Range("B9").Value  = "Group_ID"
$ws.Range("B11").Value = "Course"
$ws.Range("B12").Value = "Time Availability"
$ws.Range("B10").Value = "Group Name"
$ws.Range("B14").Value = "Admin UID"
$ws.Range("B13").Value = "Users UID"
$ws.Range("A2").Value  = "Table"
$ws.Range("B16").Value = "Chat"
$ws.Range("A3").Value  = "Columns"
$ws.Range("B18").Value = "User who posted it"
$ws.Range("B19").Value = "Time Stamp"
$ws.Range("B20").Value = "Message Contents"
$ws.Range("A8").Value  = "Table"
$ws.Range("A9").Value  = "Columns"
$ws.Range("B17").Value = "Group_ID"
$ws.Range("A16").Value = "Table"
$ws.Range("A17").Value = "Columns"

# --- Formatting --------------------------------------------------------
# Column A labels ("Table" / "Columns") are bold; column B values wrap.
# Applied top-to-bottom / left-to-right so the cellXfs records are built
# up in the same order as the target workbook (bold, then wrap, then
# bold+wrap for the "Chat" header).
$ws.Range("A2").Font.Bold  = $true
$ws.Range("B2").WrapText   = $true
$ws.Range("A3").Font.Bold  = $true
$ws.Range("B3").WrapText   = $true
$ws.Range("B4").WrapText   = $true
$ws.Range("B5").WrapText   = $true
$ws.Range("B6").WrapText   = $true
$ws.Range("A8").Font.Bold  = $true
$ws.Range("B8").WrapText   = $true
$ws.Range("A9").Font.Bold  = $true
$ws.Range("B9").WrapText   = $true
$ws.Range("B10").WrapText  = $true
$ws.Range("B11").WrapText  = $true
$ws.Range("B12").WrapText  = $true
$ws.Range("B13").WrapText  = $true
$ws.Range("B14").WrapText  = $true
$ws.Range("A16").Font.Bold = $true
$ws.Range("B16").Font.Bold = $true
$ws.Range("B16").WrapText  = $true
$ws.Range("A17").Font.Bold = $true
$ws.Range("B17").WrapText  = $true
$ws.Range("B18").WrapText  = $true
$ws.Range("B19").WrapText  = $true
$ws.Range("B20").WrapText  = $true

# --- Column widths -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.307291666666668
$ws.Columns.Item(2).ColumnWidth = 21.022135416666668
$ws.Columns.Item(3).ColumnWidth = 12.592447916666666

# --- Page setup ------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Selection ---------------------------------------------------------
$null = $ws.Range("C13").Select()
